$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 46070

$ws.Range("C3").Value = 46070

$ws.Range("C4").Value = 46070

$ws.Range("C5").Value = 46070

$ws.Range("C6").Value = 46070

$ws.Range("C7").Value = 46070

$ws.Range("C8").Value = 46070

$ws.Range("C9").Value = 46070

$ws.Range("C10").Value = 46070

$ws.Range("C11").Value = 46070

$ws.Range("C12").Value = 46070

$ws.Range("C13").Value = 46070

$ws.Range("A14").Value = "A 21421-2021"
$ws.Range("B14").Value = 44316
$ws.Range("C14").Value = 46070
$ws.Range("G14").Value = 0.6

$ws.Range("A15").Value = "A 17491-2024"
$ws.Range("B15").Value = 45415.50266203703
$ws.Range("C15").Value = 46070
$ws.Range("G15").Value = 6.2

$ws.Range("A16").Value = "A 50864-2022"
$ws.Range("B16").Value = 44867.56143518518
$ws.Range("C16").Value = 46070
$ws.Range("G16").Value = 3.3

$ws.Range("A17").Value = "A 23503-2025"
$ws.Range("B17").Value = 45795
$ws.Range("C17").Value = 46070
$ws.Range("G17").Value = 14.1

$ws.Range("A18").Value = "A 270-2025"
$ws.Range("B18").Value = 45660.48087962963
$ws.Range("C18").Value = 46070
$ws.Range("G18").Value = 8.9

$ws.Range("A19").Value = "A 49633-2024"
$ws.Range("B19").Value = 45596.59559027778
$ws.Range("C19").Value = 46070
$ws.Range("G19").Value = 0.8

$ws.Range("A20").Value = "A 24086-2025"
$ws.Range("B20").Value = 45795
$ws.Range("C20").Value = 46070
$ws.Range("G20").Value = 0.7

$ws.Range("A21").Value = "A 4422-2024"
$ws.Range("B21").Value = 45327.45375
$ws.Range("C21").Value = 46070
$ws.Range("G21").Value = 4.5

$ws.Range("A22").Value = "A 24212-2023"
$ws.Range("B22").Value = 45076
$ws.Range("C22").Value = 46070
$ws.Range("G22").Value = 5.8

$ws.Range("A23").Value = "A 51434-2025"
$ws.Range("B23").Value = 45949
$ws.Range("C23").Value = 46070
$ws.Range("G23").Value = 2.8

$ws.Range("A24").Value = "A 50239-2022"
$ws.Range("B24").Value = 44865
$ws.Range("C24").Value = 46070
$ws.Range("G24").Value = 13.2

$ws.Range("A25").Value = "A 52965-2025"
$ws.Range("B25").Value = 45956
$ws.Range("C25").Value = 46070
$ws.Range("G25").Value = 0.6

$ws.Range("A26").Value = "A 46579-2024"
$ws.Range("B26").Value = 45582.75018518518
$ws.Range("C26").Value = 46070
$ws.Range("G26").Value = 3

$ws.Range("C27").Value = 46070

$ws.Range("C28").Value = 46070

$ws.Range("A29").Value = "A 28409-2024"
$ws.Range("B29").Value = 45477.62280092593
$ws.Range("C29").Value = 46070
$ws.Range("G29").Value = 0.4

$ws.Range("C30").Value = 46070

$ws.Range("A31").Value = "A 4780-2022"
$ws.Range("B31").Value = 44592.62657407407
$ws.Range("C31").Value = 46070
$ws.Range("G31").Value = 0.9

$ws.Range("A32").Value = "A 17492-2024"
$ws.Range("B32").Value = 45415.50709490741
$ws.Range("C32").Value = 46070
$ws.Range("G32").Value = 7.7

$ws.Range("A33").Value = "A 46587-2024"
$ws.Range("B33").Value = 45582.76763888889
$ws.Range("C33").Value = 46070
$ws.Range("G33").Value = 2.8

$ws.Range("A34").Value = "A 46588-2024"
$ws.Range("B34").Value = 45582.77137731481
$ws.Range("C34").Value = 46070
$ws.Range("G34").Value = 1

$ws.Range("A35").Value = "A 28418-2024"
$ws.Range("B35").Value = 45477.62978009259
$ws.Range("C35").Value = 46070
$ws.Range("G35").Value = 0.2

$ws.Range("A36").Value = "A 49634-2024"
$ws.Range("B36").Value = 45596.59591435185
$ws.Range("C36").Value = 46070
$ws.Range("G36").Value = 0.7
